$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.705.53"
$ws.Range("E2").Value = "  -3.84%  "

$ws.Range("D3").Value = "2.910.39"
$ws.Range("E3").Value = "  -3.74%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.54%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.555"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.60%  "

$ws.Range("D9").Value = "2.913.47"
$ws.Range("E9").Value = "  -3.97%  "

$ws.Range("E10").Value = "  -4.72%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.87"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.59%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.353"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.91%  "

$ws.Range("D13").Value = "3.416.79"
$ws.Range("E13").Value = "  -3.77%  "

$ws.Range("E14").Value = "  +0.94%  "

$ws.Range("D15").Value = "60.819.66"
$ws.Range("E15").Value = "  -3.83%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.97%  "

$ws.Range("D17").Value = "2.913.33"
$ws.Range("E17").Value = "  -3.82%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000141"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.98%  "

$ws.Range("E19").Value = "  -3.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "361.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.99%  "

$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.85%  "

$ws.Range("D26").Value = "3.036.01"
$ws.Range("E26").Value = "  -3.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.452"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.182"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.32%  "

$ws.Range("E29").Value = "  +0.32%  "

$ws.Range("D30").Value = "0.0₃0862"
$ws.Range("E30").Value = "  -12.59%  "

$ws.Range("E31").Value = "  -11.84%  "

$ws.Range("E32").Value = "  +0.08%  "

$ws.Range("E33").Value = "  -4.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "154.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.08%  "

$ws.Range("E36").Value = "  -8.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.60"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.85%  "

$ws.Range("E38").Value = "  -9.19%  "

$ws.Range("E39").Value = "  -7.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.91%  "

$ws.Range("D42").Value = "2.333.48"
$ws.Range("E42").Value = "  -8.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.644"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0571"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.53%  "

$ws.Range("E49").Value = "  -6.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0926"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.64%  "
